# Update template comments: shift the values of columns J:W (rows 1-5) one
# column to the left, wrapping the J value around to the end (column W).
# This reflects the reordering of "Comment" (and its related rows) to the
# end of the block in the shared strings / header layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

for ($row = 1; $row -le 5; $row++) {
    $values = @()
    foreach ($c in $cols) {
        $values += , $ws.Range("$c$row").Value()
    }

    $first = $values[0]
    for ($i = 0; $i -lt ($cols.Length - 1); $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i + 1]
    }
    $ws.Range("$($cols[$cols.Length - 1])$row").Value = $first
}
